$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new weekly record for the first rolling window (rows 92-165) ---
# This shifts old rows 92-165 down to 93-166, dropping the previous row-166
# record's predecessor (old row 165) out of the fixed-size window below.
$ws.Rows("92").Insert()

$ws.Range("A92").Value = 10
$ws.Range("B92").Value = "Vega Modelo de Temuco"
$ws.Range("C92").Value = "La Araucanía"
$ws.Range("D92").Value = 44614
$ws.Range("E92").Value = 9
$ws.Range("F92").Value = "Fruta"
$ws.Range("G92").Value = 100108
$ws.Range("H92").Value = "Tropicales y subtropicales"
$ws.Range("I92").Value = 100108002
$ws.Range("J92").Value = "Mango"
$ws.Range("K92").Value = "Sin especificar"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 300
$ws.Range("N92").Value = 8000
$ws.Range("O92").Value = 8000
$ws.Range("P92").Value = 8000
$ws.Range("Q92").Value = "$/bandeja 4 kilos"
$ws.Range("R92").Value = "Perú"
$ws.Range("S92").Value = 2000
$ws.Range("T92").Value = 4

# Drop the oldest record of that rolling window (now duplicated at row 166)
# so the window keeps its original size.
$ws.Rows("166").Delete()

# --- Append the new weekly record for the second, still-growing list ---
# (old rows 285-303 shift down to 286-304; nothing is dropped here).
$ws.Rows("285").Insert()

$ws.Range("A285").Value = 10
$ws.Range("B285").Value = "Vega Modelo de Temuco"
$ws.Range("C285").Value = "La Araucanía"
$ws.Range("D285").Value = 44615
$ws.Range("E285").Value = 9
$ws.Range("F285").Value = "Fruta"
$ws.Range("G285").Value = 100108
$ws.Range("H285").Value = "Tropicales y subtropicales"
$ws.Range("I285").Value = 100108002
$ws.Range("J285").Value = "Mango"
$ws.Range("K285").Value = "Sin especificar"
$ws.Range("L285").Value = "Primera"
$ws.Range("M285").Value = 700
$ws.Range("N285").Value = 7000
$ws.Range("O285").Value = 7500
$ws.Range("P285").Value = 7214
$ws.Range("Q285").Value = "$/bandeja 4 kilos"
$ws.Range("R285").Value = "Perú"
$ws.Range("S285").Value = 1804
$ws.Range("T285").Value = 4
